$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'303.26"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = "'31.78"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'7.59%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.222"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'2.73%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.07357"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'8.99%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'7.836"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'6.42%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'3.733"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'8.34%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'1.492"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'6.45%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.9080"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'-0.81%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.01675"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'2,488.93%"
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'5.12%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.07496"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'8.15%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.07989"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'4.30%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.02957"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'0.96%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.09917"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Value = "'0.001490"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-5.12%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'0.04533"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'1.20%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'0.006257"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'1.61%"
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'0.51%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'2.230"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'0.03%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D22").Value = "'0.1323"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'1.23%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'4.531"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'10.30%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.1617"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'2.13%"
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'1.76%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.004425"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'7.14%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'0.0001298"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Value = "'0.0001737"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'7.31%"
$ws.Range("E28").Style = "Normal"
$ws.Range("D40").Value = "'0.04495"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'5.57%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.007202"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'5.68%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.1344"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'8.31%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.002327"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'4.31%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.01281"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'-1.56%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.00006052"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'5.40%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'1.892"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'-3.57%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.01297"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'-13.97%"
$ws.Range("E47").Style = "Normal"
